$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: Test No. 4 ---
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "7/18/2018"
$ws.Range("C5").Value = "135.75m"
$ws.Range("D5").Value = "Obstacle environment"
$ws.Range("E5").Value = 144
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0

# --- Row 6: Test No. 5 ---
$ws.Range("A4:G4").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "7/18/2018"
$ws.Range("C6").Value = "73.42m"
$ws.Range("D6").Value = "Obstacle environment"
$ws.Range("E6").Value = 81
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0

# --- Row 7: trailing blank formatted row (no border) ---
$ws.Range("A2").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7:G7").ClearContents()
$ws.Range("A7:G7").Borders.LineStyle = -4142  # xlLineStyleNone

# Update the active selection to match the saved view state
$ws.Range("F17").Select() | Out-Null

$ws.Application.CutCopyMode = $false
